{"js": "// Replace the 25 two-digit multiplication problems in the document's\n// table cells with their updated values, matching the target diff.\n// Each entry is [oldText, newText]; old values are unique in the doc,\n// so a plain body.search + Replace is safe and format-preserving.\nconst replacements = [\n  [\"63\u00d730=1890\", \"27\u00d782=2214\"],\n  [\"37\u00d785=3145\", \"74\u00d736=2664\"],\n  [\"95\u00d789=8455\", \"69\u00d761=4209\"],\n  [\"34\u00d736=1224\", \"25\u00d770=1750\"],\n  [\"52\u00d734=1768\", \"26\u00d760=1560\"],\n  [\"51\u00d799=5049\", \"71\u00d760=4260\"],\n  [\"25\u00d775=1875\", \"33\u00d736=1188\"],\n  [\"21\u00d799=2079\", \"97\u00d744=4268\"],\n  [\"37\u00d746=1702\", \"31\u00d723=713\"],\n  [\"70\u00d776=5320\", \"61\u00d755=3355\"],\n  [\"16\u00d742=672\", \"44\u00d718=792\"],\n  [\"94\u00d747=4418\", \"68\u00d763=4284\"],\n  [\"55\u00d765=3575\", \"80\u00d746=3680\"],\n  [\"85\u00d741=3485\", \"98\u00d739=3822\"],\n  [\"81\u00d798=7938\", \"36\u00d730=1080\"],\n  [\"64\u00d772=4608\", \"99\u00d775=7425\"],\n  [\"97\u00d793=9021\", \"12\u00d737=444\"],\n  [\"33\u00d747=1551\", \"72\u00d749=3528\"],\n  [\"33\u00d771=2343\", \"43\u00d778=3354\"],\n  [\"70\u00d763=4410\", \"81\u00d785=6885\"],\n  [\"43\u00d779=3397\", \"59\u00d745=2655\"],\n  [\"51\u00d758=2958\", \"99\u00d763=6237\"],\n  [\"22\u00d799=2178\", \"45\u00d797=4365\"],\n  [\"86\u00d790=7740\", \"58\u00d727=1566\"],\n  [\"86\u00d734=2924\", \"36\u00d777=2772\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit multiplication problems in the document's\n# table cells with their updated values, matching the target diff.\n# Old values are unique across the document, so Find/Replace (scoped to\n# the whole document content) is safe and format-preserving.\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{ Old = \"63\u00d730=1890\"; New = \"27\u00d782=2214\" },\n    @{ Old = \"37\u00d785=3145\"; New = \"74\u00d736=2664\" },\n    @{ Old = \"95\u00d789=8455\"; New = \"69\u00d761=4209\" },\n    @{ Old = \"34\u00d736=1224\"; New = \"25\u00d770=1750\" },\n    @{ Old = \"52\u00d734=1768\"; New = \"26\u00d760=1560\" },\n    @{ Old = \"51\u00d799=5049\"; New = \"71\u00d760=4260\" },\n    @{ Old = \"25\u00d775=1875\"; New = \"33\u00d736=1188\" },\n    @{ Old = \"21\u00d799=2079\"; New = \"97\u00d744=4268\" },\n    @{ Old = \"37\u00d746=1702\"; New = \"31\u00d723=713\" },\n    @{ Old = \"70\u00d776=5320\"; New = \"61\u00d755=3355\" },\n    @{ Old = \"16\u00d742=672\";  New = \"44\u00d718=792\" },\n    @{ Old = \"94\u00d747=4418\"; New = \"68\u00d763=4284\" },\n    @{ Old = \"55\u00d765=3575\"; New = \"80\u00d746=3680\" },\n    @{ Old = \"85\u00d741=3485\"; New = \"98\u00d739=3822\" },\n    @{ Old = \"81\u00d798=7938\"; New = \"36\u00d730=1080\" },\n    @{ Old = \"64\u00d772=4608\"; New = \"99\u00d775=7425\" },\n    @{ Old = \"97\u00d793=9021\"; New = \"12\u00d737=444\" },\n    @{ Old = \"33\u00d747=1551\"; New = \"72\u00d749=3528\" },\n    @{ Old = \"33\u00d771=2343\"; New = \"43\u00d778=3354\" },\n    @{ Old = \"70\u00d763=4410\"; New = \"81\u00d785=6885\" },\n    @{ Old = \"43\u00d779=3397\"; New = \"59\u00d745=2655\" },\n    @{ Old = \"51\u00d758=2958\"; New = \"99\u00d763=6237\" },\n    @{ Old = \"22\u00d799=2178\"; New = \"45\u00d797=4365\" },\n    @{ Old = \"86\u00d790=7740\"; New = \"58\u00d727=1566\" },\n    @{ Old = \"86\u00d734=2924\"; New = \"36\u00d777=2772\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair.New, $wdReplaceAll)\n}\n"}
